$wb = $excel.ActiveWorkbook

$wsNew = $wb.Worksheets.Item("forNewCustomer")
$wsNew.Range("C2").Value = "9881012056"
$wsNew.Range("C3").Value = "9881012057"
$wsNew.Range("C4").Value = "9881012058"
$wsNew.Range("C5").Value = "9881012059"
$wsNew.Range("C6").Value = "9881012060"
$wsNew.Range("E2").Value = "testuser156@mail.com"
$wsNew.Range("E3").Value = "testuser157@mail.com"
$wsNew.Range("E4").Value = "testuser158@mail.com"
$wsNew.Range("E5").Value = "testuser159@mail.com"
$wsNew.Range("E6").Value = "testuser160@mail.com"

$wsSync = $wb.Worksheets.Item("forSync")
$wsSync.Range("G20").Select()

$wsSearch = $wb.Worksheets.Item("searchInput")
$wsSearch.Activate()
